$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 4) + a lone tagged blank cell in row 5 (K5), matching the
# ledger template's empty-data placeholder row.
$headers = @("Date/Time","Code","Subcode","Symbol","Buy/Sell","Open/Close","Quantity","Price","Fees","Amount","Description","Account")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $headers[$i]
}
$ws.Cells.Item(5, 11).Value = " "

# Turn the header block into a real Excel Table ("Table1") spanning down to
# row 688 so new ledger rows can be appended and stay inside the table.
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A4:L688"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium2"

# Column widths tuned to fit each header label.
$ws.Columns.Item(1).ColumnWidth = 12.166666666666668
$ws.Columns.Item(2).ColumnWidth = 7.166666666666666
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 9.166666666666668
$ws.Columns.Item(5).ColumnWidth = 10.0
$ws.Columns.Item(6).ColumnWidth = 13.0
$ws.Columns.Item(7).ColumnWidth = 10.5
$ws.Columns.Item(8).ColumnWidth = 7.0
$ws.Columns.Item(9).ColumnWidth = 6.666666666666666
$ws.Columns.Item(10).ColumnWidth = 9.833333333333332
$ws.Columns.Item(11).ColumnWidth = 12.666666666666668
$ws.Columns.Item(12).ColumnWidth = 9.833333333333332

# Restore the selection that was active when the sheet was last saved.
$ws.Range("P15").Select()
